$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.747.88"
$ws.Range("E2").Value = "  +3.42%  "
$ws.Range("D3").Value = "2.638.64"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0817"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "3.055.29"
$ws.Range("E15").Value = "  +5.59%  "
$ws.Range("D16").Value = "2.650.94"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.864"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "49.645.20"
$ws.Range("E18").Value = "  +3.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "0.0₃0949"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0801"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "2.059.75"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("E47").Value = "  +13.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.42%  "
